$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Ref, $Val)
    $cell = $Sheet.Range($Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '45.740.92'
$ws.Range('E2').Value = '  -2.03%  '
Set-TextValue $ws 'D3' '2.473.02'
$ws.Range('E3').Value = '  +9.80%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue $ws 'D5' '293.72'
$ws.Range('E5').Value = '  -0.98%  '
Set-TextValue $ws 'D6' '95.19'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue $ws 'D9' '0.523'
$ws.Range('E9').Value = '  +5.09%  '
Set-TextValue $ws 'D10' '35.20'
$ws.Range('E10').Value = '  +2.80%  '
Set-TextValue $ws 'D11' '0.0785'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('E12').Value = '  +5.68%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D13' '2.851.22'
$ws.Range('E13').Value = '  +9.57%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws 'D14' '0.104'
$ws.Range('E14').Value = '  +1.92%  '
Set-TextValue $ws 'D15' '2.469.86'
$ws.Range('E15').Value = '  +9.37%  '
Set-TextValue $ws 'D16' '0.852'
$ws.Range('E16').Value = '  +8.65%  '
$ws.Range('E17').Value = '  +6.17%  '
Set-TextValue $ws 'D18' '45.827.66'
$ws.Range('E18').Value = '  -1.80%  '
Set-TextValue $ws 'D19' '12.74'
$ws.Range('E19').Value = '  +4.22%  '
Set-TextValue $ws 'D20' '0.0₃0947'
$ws.Range('E20').Value = '  -1.26%  '
Set-TextValue $ws 'D21' '6.31'
$ws.Range('E21').Value = '  +10.26%  '
Set-TextValue $ws 'D22' '67.50'
$ws.Range('E22').Value = '  +3.51%  '
Set-TextValue $ws 'D23' '245.79'
$ws.Range('E23').Value = '  +1.44%  '
Set-TextValue $ws 'D24' '2.80'
$ws.Range('E24').Value = '  +1.78%  '
Set-TextValue $ws 'D25' '1.95'
$ws.Range('E25').Value = '  +6.91%  '
$ws.Range('E26').Value = '  -0.05%  '
Set-TextValue $ws 'D27' '38.94'
$ws.Range('E27').Value = '  -3.70%  '
$ws.Range('E28').Value = '  +0.32%  '
Set-TextValue $ws 'D29' '9.82'
$ws.Range('E29').Value = '  +4.45%  '
Set-TextValue $ws 'D30' '21.90'
$ws.Range('E30').Value = '  +10.52%  '
Set-TextValue $ws 'D31' '3.80'
$ws.Range('E31').Value = '  +16.65%  '
$ws.Range('E32').Value = '  -1.69%  '
Set-TextValue $ws 'D33' '5.55'
$ws.Range('E33').Value = '  +6.02%  '
Set-TextValue $ws 'D34' '147.56'
$ws.Range('E34').Value = '  +3.34%  '
Set-TextValue $ws 'D35' '2.07'
$ws.Range('E35').Value = '  +26.74%  '
Set-TextValue $ws 'D36' '0.0773'
$ws.Range('E36').Value = '  +2.33%  '
Set-TextValue $ws 'D37' '0.114'
$ws.Range('E37').Value = '  +4.10%  '
$ws.Range('E38').Value = '  +1.18%  '
Set-TextValue $ws 'D39' '15.21'
$ws.Range('E39').Value = '  +0.77%  '
Set-TextValue $ws 'D40' '3.95'
$ws.Range('E40').Value = '  +5.35%  '
$ws.Range('E41').Value = '  +3.07%  '
Set-TextValue $ws 'D42' '2.014.84'
$ws.Range('E42').Value = '  +13.28%  '
Set-TextValue $ws 'D43' '3.24'
$ws.Range('E43').Value = '  +6.67%  '
$ws.Range('E44').Value = '  +0.16%  '
Set-TextValue $ws 'D45' '91.38'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D46' '16.40'
$ws.Range('E46').Value = '  +34.79%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D47' '1.76'
$ws.Range('E47').Value = '  -4.47%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D48' '8.63'
$ws.Range('E48').Value = '  +11.70%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D49' '103.37'
$ws.Range('E49').Value = '  +11.32%  '
Set-TextValue $ws 'D50' '2.718.30'
$ws.Range('E50').Value = '  +9.64%  '
$ws.Range('E51').Value = '  +2.79%  '
